$d = $word.ActiveDocument

# 1. Remove " / unpause" from "P - Pause / unpause"
$d.Content.Find.Execute(" / unpause", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2. Change "U - Enable DEBUG for level 2" to "U - Unpause"
#    Replace "Enable DEBUG" -> "U" (split into its own run via a formatting touch)
$full = $d.Content.Text
$idx = $full.IndexOf("Enable DEBUG")
$len = "Enable DEBUG".Length
$r1 = $d.Range($idx, $idx + $len)
$r1.Text = "U"
$r1.Bold = 1
$r1.Bold = 0

#    Replace " for level 2" -> "npause" (split into its own run via a formatting touch)
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf(" for level 2")
$len2 = " for level 2".Length
$r2 = $d.Range($idx2, $idx2 + $len2)
$r2.Text = "npause"
$r2.Bold = 1
$r2.Bold = 0
